$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole review row for zaittomer@gmail.com (row 2); this shifts
# every subsequent row up by one, exactly like the author's edit.
$ws.Rows(2).Delete()

# The row delete doesn't renumber the worksheet's hyperlink references, so
# rebuild them to match the new row numbers / relationship ids.
$ws.Hyperlinks.Delete()

$links = @(
    @("C2", "rontiddler560@gmail.com"),
    @("C3", "gregneri12@gmail.com"),
    @("C4", "snizzvered@gmail.com"),
    @("C5", "budoyoni2@gmail.com"),
    @("C7", "hermanliran@gmail.com"),
    @("C8", "gazittalia1@gmail.com"),
    @("D8", "hermanliran@gmail.com"),
    @("C9", "leviadlevi22@gmail.com"),
    @("D9", "gazittalia1@gmail.com"),
    @("C10", "freelancernachus@gmail.com"),
    @("C11", "nevilgreen@gmail.com"),
    @("D11", "vikicrestina@gmail.com"),
    @("C12", "veredsnir12@gmail.com"),
    @("D12", "kevinkors122@gmail.com"),
    @("C13", "stevewonder3001@gmail.com"),
    @("D13", "budoyoni@gmail.com"),
    @("C14", "stclerari834@gmail.com"),
    @("C15", "stcydouel274@gmail.com"),
    @("C16", "kevinkors122@gmail.com"),
    @("D16", "sinuspai@gmail.com")
)

foreach ($link in $links) {
    $ref = $link[0]
    $email = $link[1]
    $ws.Hyperlinks.Add($ws.Range($ref), "mailto:$email", "", "", $email)
}
